$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'Vtn'
$ws.Cells.Item(2,3).Value = 'Plaur'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.337313
$ws.Cells.Item(2,8).Value = 10.011939
$ws.Cells.Item(2,9).Value = 0.1958858017947999
$ws.Cells.Item(2,10).Value = 0.1958858017947999
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 14.89002333333333
$ws.Cells.Item(2,14).Value = 44.67007
$ws.Cells.Item(2,15).Value = 0.1194491234330596
$ws.Cells.Item(2,16).Value = 0.1194491234330597
$ws.Cells.Item(2,17).Value = 49.69266844063667
$ws.Cells.Item(2,18).Value = 447.23401596573
$ws.Cells.Item(2,19).Value = 0.02339838731737091
$ws.Cells.Item(2,20).Value = 0.02339838731737091

$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'Vtn'
$ws.Cells.Item(3,3).Value = 'Plaur'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.337313
$ws.Cells.Item(3,8).Value = 10.011939
$ws.Cells.Item(3,9).Value = 0.1958858017947999
$ws.Cells.Item(3,10).Value = 0.1958858017947999
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.417914
$ws.Cells.Item(3,14).Value = 19.253742
$ws.Cells.Item(3,15).Value = 0.05148509068166413
$ws.Cells.Item(3,16).Value = 0.05148509068166414
$ws.Cells.Item(3,17).Value = 21.418587825082
$ws.Cells.Item(3,18).Value = 192.767290425738
$ws.Cells.Item(3,19).Value = 0.01008519826865576
$ws.Cells.Item(3,20).Value = 0.01008519826865576

$ws.Cells.Item(4,1).Value = 'ECs'
$ws.Cells.Item(4,2).Value = 'Vtn'
$ws.Cells.Item(4,3).Value = 'Plaur'
$ws.Cells.Item(4,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.337313
$ws.Cells.Item(4,8).Value = 10.011939
$ws.Cells.Item(4,9).Value = 0.1958858017947999
$ws.Cells.Item(4,10).Value = 0.1958858017947999
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 66.36284166666667
$ws.Cells.Item(4,14).Value = 199.088525
$ws.Cells.Item(4,15).Value = 0.5323687604884161
$ws.Cells.Item(4,16).Value = 0.5323687604884162
$ws.Cells.Item(4,17).Value = 221.4735742111083
$ws.Cells.Item(4,18).Value = 1993.262167899975
$ws.Cells.Item(4,19).Value = 0.1042834814987772
$ws.Cells.Item(4,20).Value = 0.1042834814987772

$ws.Cells.Item(5,1).Value = 'ECs'
$ws.Cells.Item(5,2).Value = 'Vtn'
$ws.Cells.Item(5,3).Value = 'Plaur'
$ws.Cells.Item(5,4).Value = 'MuSCs'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.337313
$ws.Cells.Item(5,8).Value = 10.011939
$ws.Cells.Item(5,9).Value = 0.1958858017947999
$ws.Cells.Item(5,10).Value = 0.1958858017947999
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.521285666666667
$ws.Cells.Item(5,14).Value = 10.563857
$ws.Cells.Item(5,15).Value = 0.02824807435318976
$ws.Cells.Item(5,16).Value = 0.02824807435318976
$ws.Cells.Item(5,17).Value = 11.75163243208033
$ws.Cells.Item(5,18).Value = 105.764691888723
$ws.Cells.Item(5,19).Value = 0.0055333966938337
$ws.Cells.Item(5,20).Value = 0.0055333966938337

$ws.Cells.Item(6,1).Value = 'ECs'
$ws.Cells.Item(6,2).Value = 'Vtn'
$ws.Cells.Item(6,3).Value = 'Plaur'
$ws.Cells.Item(6,4).Value = 'Resolving-Mac'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.337313
$ws.Cells.Item(6,8).Value = 10.011939
$ws.Cells.Item(6,9).Value = 0.1958858017947999
$ws.Cells.Item(6,10).Value = 0.1958858017947999
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 33.46371266666667
$ws.Cells.Item(6,14).Value = 100.391138
$ws.Cells.Item(6,15).Value = 0.2684489510436703
$ws.Cells.Item(6,16).Value = 0.2684489510436703
$ws.Cells.Item(6,17).Value = 111.6788833107313
$ws.Cells.Item(6,18).Value = 1005.109949796582
$ws.Cells.Item(6,19).Value = 0.05258533801616234
$ws.Cells.Item(6,20).Value = 0.05258533801616234

$ws.Cells.Item(7,1).Value = 'FAPs'
$ws.Cells.Item(7,2).Value = 'Vtn'
$ws.Cells.Item(7,3).Value = 'Plaur'
$ws.Cells.Item(7,4).Value = 'ECs'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 10.28369633333333
$ws.Cells.Item(7,8).Value = 30.851089
$ws.Cells.Item(7,9).Value = 0.6036083824529627
$ws.Cells.Item(7,10).Value = 0.6036083824529627
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 14.89002333333333
$ws.Cells.Item(7,14).Value = 44.67007
$ws.Cells.Item(7,15).Value = 0.1194491234330596
$ws.Cells.Item(7,16).Value = 0.1194491234330597
$ws.Cells.Item(7,17).Value = 153.1244783562478
$ws.Cells.Item(7,18).Value = 1378.12030520623
$ws.Cells.Item(7,19).Value = 0.07210049218085342
$ws.Cells.Item(7,20).Value = 0.07210049218085342

$ws.Cells.Item(8,1).Value = 'FAPs'
$ws.Cells.Item(8,2).Value = 'Vtn'
$ws.Cells.Item(8,3).Value = 'Plaur'
$ws.Cells.Item(8,4).Value = 'FAPs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 10.28369633333333
$ws.Cells.Item(8,8).Value = 30.851089
$ws.Cells.Item(8,9).Value = 0.6036083824529627
$ws.Cells.Item(8,10).Value = 0.6036083824529627
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 6.417914
$ws.Cells.Item(8,14).Value = 19.253742
$ws.Cells.Item(8,15).Value = 0.05148509068166413
$ws.Cells.Item(8,16).Value = 0.05148509068166414
$ws.Cells.Item(8,17).Value = 65.99987866944866
$ws.Cells.Item(8,18).Value = 593.998908025038
$ws.Cells.Item(8,19).Value = 0.03107683230680339
$ws.Cells.Item(8,20).Value = 0.03107683230680339

$ws.Cells.Item(9,1).Value = 'FAPs'
$ws.Cells.Item(9,2).Value = 'Vtn'
$ws.Cells.Item(9,3).Value = 'Plaur'
$ws.Cells.Item(9,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 10.28369633333333
$ws.Cells.Item(9,8).Value = 30.851089
$ws.Cells.Item(9,9).Value = 0.6036083824529627
$ws.Cells.Item(9,10).Value = 0.6036083824529627
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 66.36284166666667
$ws.Cells.Item(9,14).Value = 199.088525
$ws.Cells.Item(9,15).Value = 0.5323687604884161
$ws.Cells.Item(9,16).Value = 0.5323687604884162
$ws.Cells.Item(9,17).Value = 682.4553115170805
$ws.Cells.Item(9,18).Value = 6142.097803653725
$ws.Cells.Item(9,19).Value = 0.3213422463869015
$ws.Cells.Item(9,20).Value = 0.3213422463869016

$ws.Cells.Item(10,1).Value = 'FAPs'
$ws.Cells.Item(10,2).Value = 'Vtn'
$ws.Cells.Item(10,3).Value = 'Plaur'
$ws.Cells.Item(10,4).Value = 'MuSCs'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 10.28369633333333
$ws.Cells.Item(10,8).Value = 30.851089
$ws.Cells.Item(10,9).Value = 0.6036083824529627
$ws.Cells.Item(10,10).Value = 0.6036083824529627
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.521285666666667
$ws.Cells.Item(10,14).Value = 10.563857
$ws.Cells.Item(10,15).Value = 0.02824807435318976
$ws.Cells.Item(10,16).Value = 0.02824807435318976
$ws.Cells.Item(10,17).Value = 36.21183249891923
$ws.Cells.Item(10,18).Value = 325.906492490273
$ws.Cells.Item(10,19).Value = 0.01705077446773989
$ws.Cells.Item(10,20).Value = 0.01705077446773989

$ws.Cells.Item(11,1).Value = 'FAPs'
$ws.Cells.Item(11,2).Value = 'Vtn'
$ws.Cells.Item(11,3).Value = 'Plaur'
$ws.Cells.Item(11,4).Value = 'Resolving-Mac'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 10.28369633333333
$ws.Cells.Item(11,8).Value = 30.851089
$ws.Cells.Item(11,9).Value = 0.6036083824529627
$ws.Cells.Item(11,10).Value = 0.6036083824529627
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 33.46371266666667
$ws.Cells.Item(11,14).Value = 100.391138
$ws.Cells.Item(11,15).Value = 0.2684489510436703
$ws.Cells.Item(11,16).Value = 0.2684489510436703
$ws.Cells.Item(11,17).Value = 344.1306592499202
$ws.Cells.Item(11,18).Value = 3097.175933249282
$ws.Cells.Item(11,19).Value = 0.1620380371106644
$ws.Cells.Item(11,20).Value = 0.1620380371106644

$ws.Cells.Item(12,1).Value = 'MuSCs'
$ws.Cells.Item(12,2).Value = 'Vtn'
$ws.Cells.Item(12,3).Value = 'Plaur'
$ws.Cells.Item(12,4).Value = 'ECs'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 3.355061
$ws.Cells.Item(12,8).Value = 10.065183
$ws.Cells.Item(12,9).Value = 0.196927532435664
$ws.Cells.Item(12,10).Value = 0.196927532435664
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 14.89002333333333
$ws.Cells.Item(12,14).Value = 44.67007
$ws.Cells.Item(12,15).Value = 0.1194491234330596
$ws.Cells.Item(12,16).Value = 0.1194491234330597
$ws.Cells.Item(12,17).Value = 49.95693657475666
$ws.Cells.Item(12,18).Value = 449.61242917281
$ws.Cells.Item(12,19).Value = 0.02352282112927549
$ws.Cells.Item(12,20).Value = 0.02352282112927549

$ws.Cells.Item(13,1).Value = 'MuSCs'
$ws.Cells.Item(13,2).Value = 'Vtn'
$ws.Cells.Item(13,3).Value = 'Plaur'
$ws.Cells.Item(13,4).Value = 'FAPs'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 3.355061
$ws.Cells.Item(13,8).Value = 10.065183
$ws.Cells.Item(13,9).Value = 0.196927532435664
$ws.Cells.Item(13,10).Value = 0.196927532435664
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 6.417914
$ws.Cells.Item(13,14).Value = 19.253742
$ws.Cells.Item(13,15).Value = 0.05148509068166413
$ws.Cells.Item(13,16).Value = 0.05148509068166414
$ws.Cells.Item(13,17).Value = 21.532492962754
$ws.Cells.Item(13,18).Value = 193.792436664786
$ws.Cells.Item(13,19).Value = 0.01013883186516652
$ws.Cells.Item(13,20).Value = 0.01013883186516652

$ws.Cells.Item(14,1).Value = 'MuSCs'
$ws.Cells.Item(14,2).Value = 'Vtn'
$ws.Cells.Item(14,3).Value = 'Plaur'
$ws.Cells.Item(14,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 3.355061
$ws.Cells.Item(14,8).Value = 10.065183
$ws.Cells.Item(14,9).Value = 0.196927532435664
$ws.Cells.Item(14,10).Value = 0.196927532435664
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 66.36284166666667
$ws.Cells.Item(14,14).Value = 199.088525
$ws.Cells.Item(14,15).Value = 0.5323687604884161
$ws.Cells.Item(14,16).Value = 0.5323687604884162
$ws.Cells.Item(14,17).Value = 222.6513819250083
$ws.Cells.Item(14,18).Value = 2003.862437325075
$ws.Cells.Item(14,19).Value = 0.1048380663488168
$ws.Cells.Item(14,20).Value = 0.1048380663488168

$ws.Cells.Item(15,1).Value = 'MuSCs'
$ws.Cells.Item(15,2).Value = 'Vtn'
$ws.Cells.Item(15,3).Value = 'Plaur'
$ws.Cells.Item(15,4).Value = 'MuSCs'
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 3.355061
$ws.Cells.Item(15,8).Value = 10.065183
$ws.Cells.Item(15,9).Value = 0.196927532435664
$ws.Cells.Item(15,10).Value = 0.196927532435664
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.521285666666667
$ws.Cells.Item(15,14).Value = 10.563857
$ws.Cells.Item(15,15).Value = 0.02824807435318976
$ws.Cells.Item(15,16).Value = 0.02824807435318976
$ws.Cells.Item(15,17).Value = 11.81412821009233
$ws.Cells.Item(15,18).Value = 106.327153890831
$ws.Cells.Item(15,19).Value = 0.005562823578432826
$ws.Cells.Item(15,20).Value = 0.005562823578432826

$ws.Cells.Item(16,1).Value = 'MuSCs'
$ws.Cells.Item(16,2).Value = 'Vtn'
$ws.Cells.Item(16,3).Value = 'Plaur'
$ws.Cells.Item(16,4).Value = 'Resolving-Mac'
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 3.355061
$ws.Cells.Item(16,8).Value = 10.065183
$ws.Cells.Item(16,9).Value = 0.196927532435664
$ws.Cells.Item(16,10).Value = 0.196927532435664
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 33.46371266666667
$ws.Cells.Item(16,14).Value = 100.391138
$ws.Cells.Item(16,15).Value = 0.2684489510436703
$ws.Cells.Item(16,16).Value = 0.2684489510436703
$ws.Cells.Item(16,17).Value = 112.2727972831393
$ws.Cells.Item(16,18).Value = 1010.455175548254
$ws.Cells.Item(16,19).Value = 0.05286498951397236
$ws.Cells.Item(16,20).Value = 0.05286498951397236

$ws.Cells.Item(17,1).Value = 'Resolving-Mac'
$ws.Cells.Item(17,2).Value = 'Vtn'
$ws.Cells.Item(17,3).Value = 'Plaur'
$ws.Cells.Item(17,4).Value = 'ECs'
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.06096333333333333
$ws.Cells.Item(17,8).Value = 0.18289
$ws.Cells.Item(17,9).Value = 0.003578283316573439
$ws.Cells.Item(17,10).Value = 0.003578283316573439
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 14.89002333333333
$ws.Cells.Item(17,14).Value = 44.67007
$ws.Cells.Item(17,15).Value = 0.1194491234330596
$ws.Cells.Item(17,16).Value = 0.1194491234330597
$ws.Cells.Item(17,17).Value = 0.9077454558111111
$ws.Cells.Item(17,18).Value = 8.169709102300001
$ws.Cells.Item(17,19).Value = 0.0004274228055598388
$ws.Cells.Item(17,20).Value = 0.0004274228055598388

$ws.Cells.Item(18,1).Value = 'Resolving-Mac'
$ws.Cells.Item(18,2).Value = 'Vtn'
$ws.Cells.Item(18,3).Value = 'Plaur'
$ws.Cells.Item(18,4).Value = 'FAPs'
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = 0.3333333333333333
$ws.Cells.Item(18,7).Value = 0.06096333333333333
$ws.Cells.Item(18,8).Value = 0.18289
$ws.Cells.Item(18,9).Value = 0.003578283316573439
$ws.Cells.Item(18,10).Value = 0.003578283316573439
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 6.417914
$ws.Cells.Item(18,14).Value = 19.253742
$ws.Cells.Item(18,15).Value = 0.05148509068166413
$ws.Cells.Item(18,16).Value = 0.05148509068166414
$ws.Cells.Item(18,17).Value = 0.3912574304866667
$ws.Cells.Item(18,18).Value = 3.52131687438
$ws.Cells.Item(18,19).Value = 0.0001842282410384694
$ws.Cells.Item(18,20).Value = 0.0001842282410384694

$ws.Cells.Item(19,1).Value = 'Resolving-Mac'
$ws.Cells.Item(19,2).Value = 'Vtn'
$ws.Cells.Item(19,3).Value = 'Plaur'
$ws.Cells.Item(19,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(19,5).Value = 1
$ws.Cells.Item(19,6).Value = 0.3333333333333333
$ws.Cells.Item(19,7).Value = 0.06096333333333333
$ws.Cells.Item(19,8).Value = 0.18289
$ws.Cells.Item(19,9).Value = 0.003578283316573439
$ws.Cells.Item(19,10).Value = 0.003578283316573439
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 66.36284166666667
$ws.Cells.Item(19,14).Value = 199.088525
$ws.Cells.Item(19,15).Value = 0.5323687604884161
$ws.Cells.Item(19,16).Value = 0.5323687604884162
$ws.Cells.Item(19,17).Value = 4.045700037472223
$ws.Cells.Item(19,18).Value = 36.41130033725
$ws.Cells.Item(19,19).Value = 0.00190496625392058
$ws.Cells.Item(19,20).Value = 0.00190496625392058

$ws.Cells.Item(20,1).Value = 'Resolving-Mac'
$ws.Cells.Item(20,2).Value = 'Vtn'
$ws.Cells.Item(20,3).Value = 'Plaur'
$ws.Cells.Item(20,4).Value = 'MuSCs'
$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = 0.3333333333333333
$ws.Cells.Item(20,7).Value = 0.06096333333333333
$ws.Cells.Item(20,8).Value = 0.18289
$ws.Cells.Item(20,9).Value = 0.003578283316573439
$ws.Cells.Item(20,10).Value = 0.003578283316573439
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 3.521285666666667
$ws.Cells.Item(20,14).Value = 10.563857
$ws.Cells.Item(20,15).Value = 0.02824807435318976
$ws.Cells.Item(20,16).Value = 0.02824807435318976
$ws.Cells.Item(20,17).Value = 0.2146693118588889
$ws.Cells.Item(20,18).Value = 1.93202380673
$ws.Cells.Item(20,19).Value = 0.000101079613183345
$ws.Cells.Item(20,20).Value = 0.0001010796131833449

$ws.Cells.Item(21,1).Value = 'Resolving-Mac'
$ws.Cells.Item(21,2).Value = 'Vtn'
$ws.Cells.Item(21,3).Value = 'Plaur'
$ws.Cells.Item(21,4).Value = 'Resolving-Mac'
$ws.Cells.Item(21,5).Value = 1
$ws.Cells.Item(21,6).Value = 0.3333333333333333
$ws.Cells.Item(21,7).Value = 0.06096333333333333
$ws.Cells.Item(21,8).Value = 0.18289
$ws.Cells.Item(21,9).Value = 0.003578283316573439
$ws.Cells.Item(21,10).Value = 0.003578283316573439
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 33.46371266666667
$ws.Cells.Item(21,14).Value = 100.391138
$ws.Cells.Item(21,15).Value = 0.2684489510436703
$ws.Cells.Item(21,16).Value = 0.2684489510436703
$ws.Cells.Item(21,17).Value = 2.040059469868889
$ws.Cells.Item(21,18).Value = 18.36053522882
$ws.Cells.Item(21,19).Value = 0.0009605864028712053
$ws.Cells.Item(21,20).Value = 0.0009605864028712052
